# "Update countries & provincias Spain"
#
# 1) Refresh the "last updated" timestamp banner in A1.
# 2) Add Guatemala to the country table, directly under Costa de Marfil
#    (row 84), pushing Cuba / Islandia / Republica de Macedonia / Estonia
#    down one row each (Lituania in row 89 stays put).
# 3) Add Mongolia to the country table, directly under Brunei (row 159),
#    pushing Mozambique down one row (Yemen in row 161 stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 05:35"

# --- 2) Guatemala inserted right after Costa de Marfil (row 83) ---------
# Row 84: Guatemala's own new figures (previously Cuba's row)
$ws.Range("A84").Value = "Guatemala"
$ws.Range("B84").Value = 1912
$ws.Range("C84").Value = 149
$ws.Range("D84").Value = 138
$ws.Range("E84").Value = 1739
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 35

# Row 85: Cuba (shifted down from row 84)
$ws.Range("A85").Value = "Cuba"
$ws.Range("B85").Value = 1872
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 1495
$ws.Range("E85").Value = 298
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 79

# Row 86: Islandia (shifted down from row 85)
$ws.Range("A86").Value = "Islandia"
$ws.Range("B86").Value = 1802
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 1786
$ws.Range("E86").Value = 6
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 10

# Row 87: Republica de Macedonia (shifted down from row 86)
$ws.Range("A87").Value = "Republica de Macedonia"
$ws.Range("B87").Value = 1792
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 1293
$ws.Range("E87").Value = 398
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 101

# Row 88: Estonia (shifted down from row 87)
$ws.Range("A88").Value = "Estonia"
$ws.Range("B88").Value = 1774
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 938
$ws.Range("E88").Value = 773
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 63

# Row 89 (Lituania) is unchanged.

# --- 3) Mongolia inserted right after Brunei (row 158) -------------------
# Row 159: Mongolia's own new figures (previously Mozambique's row)
$ws.Range("A159").Value = "Mongolia"
$ws.Range("B159").Value = 140
$ws.Range("C159").Value = 4
$ws.Range("D159").Value = 24
$ws.Range("E159").Value = 116
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0

# Row 160: Mozambique (shifted down from row 159)
$ws.Range("A160").Value = "Mozambique"
$ws.Range("B160").Value = 137
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 44
$ws.Range("E160").Value = 93
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0

# Row 161 (Yemen) is unchanged.
